$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 807.5
$ws.Range("J17").Value = 848.82355
$ws.Range("L17").Value = 2546.47065
$ws.Range("N17").Value = -2882.47065

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 313.63635
$ws.Range("I107").Value = 312
$ws.Range("J107").Value = 315.6
$ws.Range("K107").Value = 312
$ws.Range("L107").Value = 315.6
$ws.Range("M107").Value = 1608
$ws.Range("N107").Value = -4155.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3691.6843
$ws.Range("I113").Value = 2505.111
$ws.Range("K113").Value = 2505.111
$ws.Range("M113").Value = 748.8890000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1795.5555
$ws.Range("I129").Value = 1270
$ws.Range("J129").Value = 6000
$ws.Range("K129").Value = 3810
$ws.Range("L129").Value = 18000
$ws.Range("M129").Value = 1190
$ws.Range("N129").Value = -28000

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 485678.1
$ws.Range("I137").Value = 1387.5
$ws.Range("K137").Value = 4162.5
$ws.Range("M137").Value = -1612.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1776.7222
$ws.Range("I138").Value = 1415.5834
$ws.Range("J138").Value = 2499
$ws.Range("K138").Value = 4246.7502
$ws.Range("L138").Value = 7497
$ws.Range("M138").Value = 893.2497999999996
$ws.Range("N138").Value = -17777

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2515.96
$ws.Range("I141").Value = 2204.125
$ws.Range("K141").Value = 6612.375
$ws.Range("M141").Value = -1432.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 33018.28
$ws.Range("I61").Value = 1879.2727
$ws.Range("J61").Value = 101524.1
$ws.Range("K61").Value = 1879.2727
$ws.Range("L61").Value = 101524.1
$ws.Range("M61").Value = -1667.2727
$ws.Range("N61").Value = -101948.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 32538.818
$ws.Range("I74").Value = 51434.7
$ws.Range("K74").Value = 51434.7
$ws.Range("M74").Value = -50560.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 32538.818
$ws.Range("I77").Value = 51434.7
$ws.Range("K77").Value = 257173.5
$ws.Range("M77").Value = -252805.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3042.32
$ws.Range("I122").Value = 1400.8182
$ws.Range("K122").Value = 4202.4546
$ws.Range("M122").Value = -1752.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2493.7932
$ws.Range("I132").Value = 2494.0527
$ws.Range("J132").Value = 2493.3
$ws.Range("K132").Value = 7482.158100000001
$ws.Range("L132").Value = 7479.900000000001
$ws.Range("M132").Value = -4952.158100000001
$ws.Range("N132").Value = -12539.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 33018.28
$ws.Range("I136").Value = 1879.2727
$ws.Range("J136").Value = 101524.1
$ws.Range("K136").Value = 5637.8181
$ws.Range("L136").Value = 304572.3
$ws.Range("M136").Value = -3087.8181
$ws.Range("N136").Value = -309672.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3899.125
$ws.Range("I94").Value = 3899.125
$ws.Range("K94").Value = 3899.125
$ws.Range("M94").Value = -3448.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H111").Value = 100702
$ws.Range("J111").Value = 100702
$ws.Range("L111").Value = 100702
$ws.Range("N111").Value = -108882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3776.5
$ws.Range("I134").Value = 2468.5833
$ws.Range("K134").Value = 7405.7499
$ws.Range("M134").Value = -4870.7499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3265.5144
$ws.Range("I31").Value = 2656.3809
$ws.Range("J31").Value = 4179.2144
$ws.Range("K31").Value = 2656.3809
$ws.Range("L31").Value = 4179.2144
$ws.Range("M31").Value = -2361.3809
$ws.Range("N31").Value = -4769.2144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3265.5144
$ws.Range("I34").Value = 2656.3809
$ws.Range("J34").Value = 4179.2144
$ws.Range("K34").Value = 2656.3809
$ws.Range("L34").Value = 4179.2144
$ws.Range("M34").Value = -2454.3809
$ws.Range("N34").Value = -4583.2144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 6347.25
$ws.Range("I39").Value = 6347.25
$ws.Range("K39").Value = 6347.25
$ws.Range("M39").Value = -5956.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 6347.25
$ws.Range("I49").Value = 6347.25
$ws.Range("K49").Value = 6347.25
$ws.Range("M49").Value = -6165.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1871.0358
$ws.Range("I58").Value = 1729.2222
$ws.Range("J58").Value = 2126.3
$ws.Range("K58").Value = 1729.2222
$ws.Range("L58").Value = 2126.3
$ws.Range("M58").Value = -1526.2222
$ws.Range("N58").Value = -2532.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7911773
$ws.Range("I99").Value = 18521096
$ws.Range("J99").Value = 2607112
$ws.Range("K99").Value = 18521096
$ws.Range("L99").Value = 2607112
$ws.Range("M99").Value = -18519598
$ws.Range("N99").Value = -2610108

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 61659.844
$ws.Range("I105").Value = 81416.92999999999
$ws.Range("K105").Value = 81416.92999999999
$ws.Range("M105").Value = -79669.92999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 7911773
$ws.Range("I126").Value = 18521096
$ws.Range("J126").Value = 2607112
$ws.Range("K126").Value = 55563288
$ws.Range("L126").Value = 7821336
$ws.Range("M126").Value = -55560818
$ws.Range("N126").Value = -7826276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1463149.2
$ws.Range("I132").Value = 1569096.1
$ws.Range("J132").Value = 1183834.8
$ws.Range("K132").Value = 4707288.300000001
$ws.Range("L132").Value = 3551504.4
$ws.Range("M132").Value = -4704758.300000001
$ws.Range("N132").Value = -3556564.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3295115
$ws.Range("I134").Value = 5497964
$ws.Range("J134").Value = 113221.445
$ws.Range("K134").Value = 16493892
$ws.Range("L134").Value = 339664.335
$ws.Range("M134").Value = -16491357
$ws.Range("N134").Value = -344734.335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1871.0358
$ws.Range("I136").Value = 1729.2222
$ws.Range("J136").Value = 2126.3
$ws.Range("K136").Value = 5187.6666
$ws.Range("L136").Value = 6378.900000000001
$ws.Range("M136").Value = -2637.6666
$ws.Range("N136").Value = -11478.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 150.76923
$ws.Range("J2").Value = 61.42857
$ws.Range("L2").Value = 368.57142
$ws.Range("N2").Value = -594.57142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 125241.375
$ws.Range("I6").Value = 125241.375
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 375724.125
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -375611.125
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 350
$ws.Range("J11").Value = 350
$ws.Range("L11").Value = 1050
$ws.Range("N11").Value = -1330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 9854.200000000001
$ws.Range("I32").Value = 20135.5
$ws.Range("K32").Value = 60406.5
$ws.Range("M32").Value = -60123.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 769.75
$ws.Range("I38").Value = 45
$ws.Range("J38").Value = 873.2857
$ws.Range("K38").Value = 135
$ws.Range("L38").Value = 2619.8571
$ws.Range("M38").Value = 212
$ws.Range("N38").Value = -3313.8571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1714.8889
$ws.Range("I140").Value = 1398
$ws.Range("K140").Value = 4194
$ws.Range("M140").Value = 986

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3167.3333
$ws.Range("I141").Value = 3167.3333
$ws.Range("K141").Value = 9501.999899999999
$ws.Range("M141").Value = -4321.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 13749.5
$ws.Range("J44").Value = 17500
$ws.Range("L44").Value = 17500
$ws.Range("N44").Value = -18692

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 842.2353000000001
$ws.Range("I97").Value = 213.27272
$ws.Range("K97").Value = 213.27272
$ws.Range("M97").Value = 282.72728

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4708.4
$ws.Range("I132").Value = 3259.9443
$ws.Range("K132").Value = 9779.832900000001
$ws.Range("M132").Value = -7249.832900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6176707.5
$ws.Range("I40").Value = 3633.8333
$ws.Range("J40").Value = 18522854
$ws.Range("K40").Value = 3633.8333
$ws.Range("L40").Value = 18522854
$ws.Range("M40").Value = -3497.8333
$ws.Range("N40").Value = -18523126

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 113368.31
$ws.Range("I68").Value = 177933.25
$ws.Range("J68").Value = 2685.5715
$ws.Range("K68").Value = 177933.25
$ws.Range("L68").Value = 2685.5715
$ws.Range("M68").Value = -177184.25
$ws.Range("N68").Value = -4183.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 113368.31
$ws.Range("I71").Value = 177933.25
$ws.Range("J71").Value = 2685.5715
$ws.Range("K71").Value = 889666.25
$ws.Range("L71").Value = 13427.8575
$ws.Range("M71").Value = -885922.25
$ws.Range("N71").Value = -20915.8575

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2481.2666
$ws.Range("I132").Value = 2776.625
$ws.Range("J132").Value = 2143.7144
$ws.Range("K132").Value = 8329.875
$ws.Range("L132").Value = 6431.1432
$ws.Range("M132").Value = -5799.875
$ws.Range("N132").Value = -11491.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2636.125
$ws.Range("I136").Value = 2025.125
$ws.Range("J136").Value = 3247.125
$ws.Range("K136").Value = 6075.375
$ws.Range("L136").Value = 9741.375
$ws.Range("M136").Value = -3525.375
$ws.Range("N136").Value = -14841.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 97714.75
$ws.Range("J139").Value = 97714.75
$ws.Range("L139").Value = 97714.75
$ws.Range("N139").Value = -107994.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1300.7778
$ws.Range("I122").Value = 1300.7778
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3902.3334
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1452.3334
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2497.8
$ws.Range("I132").Value = 2372.923
$ws.Range("J132").Value = 2729.7144
$ws.Range("K132").Value = 7118.768999999999
$ws.Range("L132").Value = 8189.1432
$ws.Range("M132").Value = -4588.768999999999
$ws.Range("N132").Value = -13249.1432
